$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "29÷5="
$t.Cell(1, 2).Range.Text = "67÷2="
$t.Cell(1, 3).Range.Text = "61÷3="
$t.Cell(1, 4).Range.Text = "18÷7="
$t.Cell(1, 5).Range.Text = "54÷6="
$t.Cell(5, 1).Range.Text = "73÷8="
$t.Cell(5, 2).Range.Text = "94÷2="
$t.Cell(5, 3).Range.Text = "87÷9="
$t.Cell(5, 4).Range.Text = "10÷8="
$t.Cell(5, 5).Range.Text = "46÷7="
$t.Cell(9, 1).Range.Text = "17÷4="
$t.Cell(9, 2).Range.Text = "79÷8="
$t.Cell(9, 3).Range.Text = "43÷6="
$t.Cell(9, 4).Range.Text = "91÷5="
$t.Cell(9, 5).Range.Text = "63÷7="
$t.Cell(13, 1).Range.Text = "68÷4="
$t.Cell(13, 2).Range.Text = "98÷4="
$t.Cell(13, 3).Range.Text = "39÷3="
$t.Cell(13, 4).Range.Text = "59÷5="
$t.Cell(13, 5).Range.Text = "32÷4="
$t.Cell(17, 1).Range.Text = "28÷9="
$t.Cell(17, 2).Range.Text = "82÷4="
$t.Cell(17, 3).Range.Text = "22÷3="
$t.Cell(17, 4).Range.Text = "15÷6="
$t.Cell(17, 5).Range.Text = "34÷4="
